$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1 "PI hours": add a new "app" column (F) that keeps the original
# list-style affiliation strings (e.g. "['ECE', 'CSL']"), and simplify the
# existing "dept" column (E) down to a single department code.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PI hours")

# Remember the original "dept" values before overwriting column E.
$origDept2 = $ws1.Range("E2").Value2
$origDept3 = $ws1.Range("E3").Value2
$origDept4 = $ws1.Range("E4").Value2

$ws1.Range("F1").Value2 = "app"
$ws1.Range("E1").Copy() | Out-Null
$ws1.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws1.Range("F2").Value2 = $origDept2
$ws1.Range("F3").Value2 = $origDept3
$ws1.Range("F4").Value2 = $origDept4

$ws1.Range("E2").Value2 = "ECE"
$ws1.Range("E3").Value2 = "ME"
$ws1.Range("E4").Value2 = "ECE"

# ---------------------------------------------------------------------------
# New sheet "unit(accumulative) hours": a full copy of the still-untouched
# "dept hours" sheet (CSL / ECE / ME), with only its header relabeled from
# "dept" to "unit(accumulative)". Do this *before* editing "dept hours" so
# the clone naturally keeps the CSL row and all original formatting.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dept hours")
$ws2.Copy($null, $ws2) | Out-Null
$ws3 = $wb.Worksheets.Item($ws2.Index + 1)
$ws3.Name = "unit(accumulative) hours"
$ws3.Range("B1").Value2 = "unit(accumulative)"

# ---------------------------------------------------------------------------
# Sheet2 "dept hours" -> renamed "department hours": drop the CSL row (it
# now lives on the new accumulative sheet) and recompute hours/percentage
# for ECE / ME only.
# ---------------------------------------------------------------------------
$ws2.Name = "department hours"
$ws2.Rows.Item(4).Delete() | Out-Null

$ws2.Range("B2").Value2 = "ECE"
$ws2.Range("C2").Value2 = 120
$ws2.Range("D2").Value2 = 73.17073170731707

$ws2.Range("B3").Value2 = "ME"
$ws2.Range("C3").Value2 = 44
$ws2.Range("D3").Value2 = 26.82926829268293

# Restore "PI hours" as the active/selected tab, matching the original file.
$ws1.Activate()
